$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.797.72"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "'2.042.29"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'227.50"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'0.608"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("D7").Value = "'59.74"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.50%  "
$ws.Range("D10").Value = "'0.0834"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "'2.342.40"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'14.47"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "'21.08"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("E15").Value = "  +5.54%  "
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "'2.041.76"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'37.772.67"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "'69.53"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'5.90"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").Value = "'0.0₃0824"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'223.88"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "'2.43"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").Value = "'168.90"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").Value = "'18.80"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +9.08%  "
$ws.Range("D33").Value = "'4.39"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "'0.0604"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").Value = "'6.51"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").Value = "'2.34"
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("D38").Value = "'3.43"
$ws.Range("E38").Value = "  +5.75%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'18.05"
$ws.Range("E40").Value = "  +6.61%  "
$ws.Range("D41").Value = "'1.526.97"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").Value = "'97.46"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "'4.22"
$ws.Range("E45").Value = "  +7.95%  "
$ws.Range("D46").Value = "'0.0905"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "'7.08"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("D51").Value = "'2.231.98"
$ws.Range("E51").Value = "  +0.42%  "
